$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: update the daily conversion text with new rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$r = $wsHoja1.Range("A1")
$t = $r.Text
$t = $t.Replace("1000 Bs = 2.33 = 8779.46 pesos", "1000 Bs = 2.35 = 8801.22 pesos")
$t = $t.Replace("8779.46 pesos = 2.32 = 958.44 Bs", "8801.22 pesos = 2.33 = 940.58 Bs")
$r.Value = $t

# --- tasas sheet: update N10/O10 and N12/O12 values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 426.1
$wsTasas.Range("O10").Value = 3750.2
$wsTasas.Range("N12").Value = 3785
$wsTasas.Range("O12").Value = 404.5
